# Adds "start_col", "end_col", "start_col_no" and "end_col_no" columns
# (C:F) to the ref-type dictionary sheet, giving the Excel column-letter
# range and numeric column-index range that corresponds to each XML
# reference type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is Row, Column letter, Value - written in the exact order
# the values were originally entered into the sheet.
$writes = @(
    @(1,  "C", "start_col"),
    @(1,  "D", "end_col"),

    @(2,  "C", "AG"),
    @(2,  "D", "AR"),
    @(3,  "C", "AS"),
    @(3,  "D", "BJ"),
    @(4,  "C", "AS"),
    @(4,  "D", "BJ"),
    @(5,  "C", "AS"),
    @(5,  "D", "BJ"),
    @(6,  "C", "DA"),
    @(6,  "D", "DR"),
    @(7,  "C", "DS"),
    @(7,  "D", "EL"),
    @(8,  "C", "BK"),
    @(8,  "D", "BV"),
    @(9,  "D", "CM"),
    @(9,  "C", "BW"),
    @(10, "C", "CN"),
    @(10, "D", "CZ"),
    @(11, "C", "O"),
    @(11, "D", "W"),
    @(12, "C", "X"),
    @(12, "D", "AF"),
    @(13, "C", "EM"),
    @(13, "D", "FA"),
    @(14, "C", "G"),
    @(14, "D", "N"),

    @(1,  "E", "start_col_no"),
    @(1,  "F", "end_col_no"),

    @(2,  "E", 33),
    @(2,  "F", 44),
    @(3,  "E", 45),
    @(3,  "F", 62),
    @(4,  "E", 45),
    @(4,  "F", 62),
    @(5,  "E", 45),
    @(5,  "F", 62),
    @(6,  "E", 105),
    @(6,  "F", 122),
    @(7,  "E", 123),
    @(7,  "F", 142),
    @(8,  "E", 63),
    @(8,  "F", 74),
    @(9,  "E", 75),
    @(9,  "F", 91),
    @(10, "E", 92),
    @(10, "F", 104),
    @(11, "E", 15),
    @(11, "F", 23),
    @(12, "E", 24),
    @(12, "F", 32),
    @(13, "E", 143),
    @(13, "F", 157),
    @(14, "E", 6),
    @(14, "F", 14)
)

foreach ($w in $writes) {
    $rowNum = $w[0]
    $colLetter = $w[1]
    $value = $w[2]
    $ws.Range("$colLetter$rowNum").Value = $value
}
